# Update "想去人数" (F column) values across sheets to reflect newly
# generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F4").Value = 353
$ws1.Range("F5").Value = 5063
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 9616
$ws1.Range("F8").Value = 0
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 0

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("F6").Value = 0

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1228
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 14
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 0
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 698
$ws4.Range("F18").Value = 0
